$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 219; this shifts existing rows 219..305 down to 220..306
$ws.Rows.Item(219).Insert()

# Populate the newly inserted row 219 with the new data record
$ws.Range("A219").Value = 4
$ws.Range("B219").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C219").Value = "Los Lagos"
$ws.Range("D219").Value = 44784
$ws.Range("E219").Value = 10
$ws.Range("F219").Value = 100112043
$ws.Range("G219").Value = "Pepino ensalada"
$ws.Range("H219").Value = "Sin especificar"
$ws.Range("I219").Value = "Primera"
$ws.Range("J219").Value = 120
$ws.Range("K219").Value = 27000
$ws.Range("L219").Value = 27000
$ws.Range("M219").Value = 27000
$ws.Range("N219").Value = "`$/caja 60 unidades"
$ws.Range("O219").Value = "Región de Arica y Parinacota"
$ws.Range("P219").Value = 450
$ws.Range("Q219").Value = 60
$ws.Range("R219").Value = "Hortaliza"
